# Update "想去人数" (want-to-go count) figures in column F across all sheets.
# "全部类型" (sheet 4) mirrors rows from the other three sheets, so matching rows are updated there too.
$wb = $excel.ActiveWorkbook

# --- 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 5925
$ws.Range("F6").Value = 64
$ws.Range("F13").Value = 1592
$ws.Range("F14").Value = 1592
$ws.Range("F15").Value = 1584
$ws.Range("F16").Value = 556
$ws.Range("F17").Value = 170
$ws.Range("F18").Value = 631
$ws.Range("F19").Value = 4510
$ws.Range("F20").Value = 34
$ws.Range("F22").Value = 3347
$ws.Range("F23").Value = 820
$ws.Range("F26").Value = 3
$ws.Range("F27").Value = 2318
$ws.Range("F33").Value = 788
$ws.Range("F37").Value = 1224
$ws.Range("F38").Value = 1204
# --- 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F20").Value = 235
# --- 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 704
$ws.Range("F4").Value = 192
$ws.Range("F5").Value = 282
# --- 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 704
$ws.Range("F6").Value = 192
$ws.Range("F7").Value = 5925
$ws.Range("F8").Value = 5925
$ws.Range("F10").Value = 64
$ws.Range("F23").Value = 1592
$ws.Range("F25").Value = 1584
$ws.Range("F26").Value = 556
$ws.Range("F27").Value = 170
$ws.Range("F28").Value = 631
$ws.Range("F29").Value = 4510
$ws.Range("F31").Value = 3347
$ws.Range("F32").Value = 820
$ws.Range("F36").Value = 2318
$ws.Range("F44").Value = 235
$ws.Range("F46").Value = 788
$ws.Range("F48").Value = 1224
$ws.Range("F50").Value = 1204
